$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header "iD" in column E
$ws.Range("E1").Value = "iD"

# Add numeric data for the new "iD" column
$ws.Range("E2").Value = 232
$ws.Range("E3").Value = 233
$ws.Range("E4").Value = 234

# Update the active selection to match the saved state
$ws.Range("E7").Select()
